$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 2
    3 = 2
    4 = 1
    5 = 1
    6 = 2
    7 = 2
    8 = 2
    9 = 2
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 0
    17 = 2
    18 = 1
    19 = 0
    20 = 2
    21 = 0
    22 = 0
    23 = 1
    24 = 1
    25 = 0
    26 = 0
    27 = 0
    28 = 2
    29 = 2
    30 = 2
    31 = 0
    32 = 1
    33 = 2
    34 = 0
    35 = 2
    36 = 0
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 1
    42 = 0
    43 = 0
    44 = 1
    45 = 1
    46 = 2
    47 = 1
    48 = 1
    49 = 0
    50 = 2
    51 = 2
    52 = 1
    53 = 1
    54 = 0
    55 = 0
    56 = 0
    57 = 2
    58 = 1
    59 = 2
    60 = 1
    61 = 2
    62 = 1
    63 = 2
    64 = 2
    65 = 0
    66 = 2
    67 = 2
    68 = 1
    69 = 0
    70 = 1
    71 = 0
    72 = 2
    73 = 1
    74 = 0
    75 = 1
    76 = 2
    77 = 1
    78 = 0
    79 = 0
    80 = 2
    81 = 1
    82 = 0
    83 = 0
    84 = 0
    85 = 2
    86 = 1
    87 = 1
    88 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}